$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 487
$ws1.Range("F8").Value = 9
$ws1.Range("F14").Value = 5930
$ws1.Range("F19").Value = 1246
$ws1.Range("F23").Value = 820
$ws1.Range("F26").Value = 10752
$ws1.Range("F29").Value = 122
$ws1.Range("F31").Value = 252

# Sheet "全部类型" (sheet4.xml) - same events, column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 487
$ws4.Range("F9").Value = 9
$ws4.Range("F15").Value = 5930
$ws4.Range("F20").Value = 1246
$ws4.Range("F24").Value = 820
$ws4.Range("F28").Value = 10752
$ws4.Range("F31").Value = 122
$ws4.Range("F33").Value = 252
